$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated error-table values (naive component forecaster bug fix)
$data = @{
    2  = @(-0.011922918044421,   1.502158260014067,  8.459308576868294,  2.908489053936475, 2.973831532800685, 23)
    3  = @( 0.5658384863776997,  2.129283520569611, 18.65469822571722,  4.319108498951747, 4.382647333937844, 22)
    4  = @(-0.9240374313429358,  1.444924223990262,  8.174846682920178,  2.85916887974813,  2.772553260771954, 21)
    5  = @( 0.07040277154089039, 0.5120019363376737, 0.4676573765475683, 0.6838547920045368, 0.6978921846357898, 20)
    6  = @( 0.1483703450672545,  0.8649345531276215, 1.90040254551332,   1.378550886080496,  1.408099328451515, 19)
    7  = @(-0.07978550902713241, 0.7076988813810415, 1.325148457270131,  1.151150927233319,  1.181676006574844, 18)
    8  = @(-0.004526189646546301,0.5864325954492303, 0.5709791096934894, 0.7556315965425807, 0.7788732485180572,17)
    9  = @( 0.2502318746956579,  0.5174934839310298, 0.4275779196144051, 0.653894425434569,   0.6239332686546548,16)
    10 = @( 0.1901465307434204,  0.4556814281592019, 0.3317531417364475, 0.575980157415555,   0.5627713046510059,15)
    11 = @( 0.2349662137414842,  0.4912372055275555, 0.3821517403275286, 0.6181842284687702,  0.5951368220458327,13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
